$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Values -------------------------------------------------------------
# Row 1 header-ish values
$ws.Range("A1").Value = 1
$ws.Range("B1").Value = 2
$ws.Range("C1").Value = 3
$ws.Range("D1").Value = 4
$ws.Range("E1").Value = 5

# Column A values for rows 2..5
$ws.Range("A2").Value = 2
$ws.Range("A3").Value = 3
$ws.Range("A4").Value = 4
$ws.Range("A5").Value = 5

# Diagonal "1" markers
$ws.Range("B2").Value = 1
$ws.Range("C3").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("E5").Value = 1

# ---- Formatting -----------------------------------------------------------
$rng = $ws.Range("A1:E5")

# Yellow fill across the whole block
$rng.Interior.Color = 65535

# Medium-weight border around the outside of the block
$xlEdgeLeft = 7
$xlEdgeTop = 8
$xlEdgeBottom = 9
$xlEdgeRight = 10
$xlMedium = -4138

$rng.Borders.Item($xlEdgeLeft).Weight = $xlMedium
$rng.Borders.Item($xlEdgeTop).Weight = $xlMedium
$rng.Borders.Item($xlEdgeRight).Weight = $xlMedium
$rng.Borders.Item($xlEdgeBottom).Weight = $xlMedium

# ---- Selection / view state ------------------------------------------------
$rng.Select()

# ---- Misc workbook metadata -------------------------------------------------
$wb.Styles.Item("Normal").Name = "Обычный"
